# edit.ps1 - applies "New crime data collected" update to the 116th Precinct weekly CompStat sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the report header text (volume number, week-covering date range) ---
$ws.Range("A8").Value = "Volume 33   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/5/2026  Through  1/11/2026"

# --- 2. Update the weekly crime-complaint statistics table (rows 15-28) ---
# Row 15
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = 0
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 1
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("G15").Value = 2
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H15").Value = -50
$ws.Range("I15").NumberFormat = "#,##0"
$ws.Range("I15").Value = 1
$ws.Range("J15").NumberFormat = "#,##0"
$ws.Range("J15").Value = 1
$ws.Range("K15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K15").Value = 0

# Row 16
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("C16").Value = 3
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("F16").NumberFormat = "#,##0"
$ws.Range("F16").Value = 6
$ws.Range("G16").NumberFormat = "#,##0"
$ws.Range("G16").Value = 6
$ws.Range("H16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H16").Value = 0
$ws.Range("I16").NumberFormat = "#,##0"
$ws.Range("I16").Value = 4
$ws.Range("K16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K16").Value = 100

# Row 17
$ws.Range("D17").NumberFormat = "#,##0"
$ws.Range("D17").Value = 6
$ws.Range("E17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E17").Value = -33.333333333333
$ws.Range("G17").NumberFormat = "#,##0"
$ws.Range("G17").Value = 25
$ws.Range("H17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H17").Value = -36
$ws.Range("I17").NumberFormat = "#,##0"
$ws.Range("I17").Value = 8
$ws.Range("J17").NumberFormat = "#,##0"
$ws.Range("J17").Value = 10
$ws.Range("K17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K17").Value = -20

# Row 18
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 1
$ws.Range("E18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E18").Value = -100
$ws.Range("F18").NumberFormat = "#,##0"
$ws.Range("F18").Value = 5
$ws.Range("H18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H18").Value = 25
$ws.Range("J18").NumberFormat = "#,##0"
$ws.Range("J18").Value = 2

# Row 19
$ws.Range("C19").NumberFormat = "#,##0"
$ws.Range("C19").Value = 4
$ws.Range("D19").NumberFormat = "#,##0"
$ws.Range("D19").Value = 7
$ws.Range("E19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E19").Value = -42.857142857142
$ws.Range("F19").NumberFormat = "#,##0"
$ws.Range("F19").Value = 15
$ws.Range("G19").NumberFormat = "#,##0"
$ws.Range("G19").Value = 25
$ws.Range("H19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H19").Value = -40
$ws.Range("I19").NumberFormat = "#,##0"
$ws.Range("I19").Value = 7
$ws.Range("J19").NumberFormat = "#,##0"
$ws.Range("J19").Value = 8
$ws.Range("K19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K19").Value = -12.5
$ws.Range("L19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L19").Value = -41.666666666666

# Row 20
$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("C20").Value = 2
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("D20").Value = 6
$ws.Range("E20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").NumberFormat = "#,##0"
$ws.Range("F20").Value = 10
$ws.Range("G20").NumberFormat = "#,##0"
$ws.Range("G20").Value = 12
$ws.Range("H20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H20").Value = -16.666666666666
$ws.Range("I20").NumberFormat = "#,##0"
$ws.Range("I20").Value = 2
$ws.Range("J20").NumberFormat = "#,##0"
$ws.Range("J20").Value = 7
$ws.Range("K20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K20").Value = -71.428571428571
$ws.Range("L20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L20").Value = -60

# Row 21
$ws.Range("C21").NumberFormat = "#,##0"
$ws.Range("C21").Value = 14
$ws.Range("D21").NumberFormat = "#,##0"
$ws.Range("D21").Value = 21
$ws.Range("E21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").NumberFormat = "#,##0"
$ws.Range("F21").Value = 53
$ws.Range("G21").NumberFormat = "#,##0"
$ws.Range("G21").Value = 74
$ws.Range("H21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("H21").Value = -28.378378378378
$ws.Range("I21").NumberFormat = "#,##0"
$ws.Range("I21").Value = 22
$ws.Range("J21").NumberFormat = "#,##0"
$ws.Range("J21").Value = 30
$ws.Range("K21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("K21").Value = -26.666666666666
$ws.Range("L21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("L21").Value = -24.137931034482

# Row 24
$ws.Range("C24").NumberFormat = "#,##0"
$ws.Range("C24").Value = 5
$ws.Range("D24").NumberFormat = "#,##0"
$ws.Range("D24").Value = 7
$ws.Range("E24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E24").Value = -28.571428571428
$ws.Range("F24").NumberFormat = "#,##0"
$ws.Range("F24").Value = 25
$ws.Range("G24").NumberFormat = "#,##0"
$ws.Range("G24").Value = 36
$ws.Range("H24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H24").Value = -30.555555555555
$ws.Range("I24").NumberFormat = "#,##0"
$ws.Range("I24").Value = 8
$ws.Range("J24").NumberFormat = "#,##0"
$ws.Range("J24").Value = 10
$ws.Range("K24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K24").Value = -20
$ws.Range("L24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L24").Value = -33.333333333333

# Row 25
$ws.Range("C25").NumberFormat = "#,##0"
$ws.Range("C25").Value = 2
$ws.Range("D25").NumberFormat = "#,##0"
$ws.Range("D25").Value = 1
$ws.Range("E25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E25").Value = 100
$ws.Range("F25").NumberFormat = "#,##0"
$ws.Range("F25").Value = 9
$ws.Range("G25").NumberFormat = "#,##0"
$ws.Range("G25").Value = 9
$ws.Range("H25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H25").Value = 0
$ws.Range("I25").NumberFormat = "#,##0"
$ws.Range("I25").Value = 3
$ws.Range("J25").NumberFormat = "#,##0"
$ws.Range("J25").Value = 3
$ws.Range("K25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K25").Value = 0
$ws.Range("L25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L25").Value = 200

# Row 26
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("C26").Value = 7
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 4
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E26").Value = 75
$ws.Range("F26").NumberFormat = "#,##0"
$ws.Range("F26").Value = 23
$ws.Range("H26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H26").Value = -25.806451612903
$ws.Range("I26").NumberFormat = "#,##0"
$ws.Range("I26").Value = 10
$ws.Range("J26").NumberFormat = "#,##0"
$ws.Range("J26").Value = 8
$ws.Range("K26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K26").Value = 25
$ws.Range("L26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L26").Value = 100

# Row 27
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = 0
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 1
$ws.Range("G27").NumberFormat = "#,##0"
$ws.Range("G27").Value = 2
$ws.Range("H27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H27").Value = -50
$ws.Range("I27").NumberFormat = "#,##0"
$ws.Range("I27").Value = 1
$ws.Range("J27").NumberFormat = "#,##0"
$ws.Range("J27").Value = 1
$ws.Range("K27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K27").Value = 0

# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "0"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "***.*"

# --- 3. Restore the shared "plain text" cell style (s13-equivalent) on cells that became text ---
# (NumberFormat = "@" above created a distinct style; copy the formatting back from a
#  known text cell so the resulting style matches the rest of the worksheet.)
$ws.Range("C14").Copy() | Out-Null
$textFixupRanges = @("D16","E16","C18","G28","H28")
foreach ($addr in $textFixupRanges) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
